$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5087719298245614
$ws.Range("B3").Value = 0.5423728813559322
$ws.Range("B4").Value = 0.5333333333333334
$ws.Range("B5").Value = 0.5245901639344263
$ws.Range("B6").Value = 43.06816697120667
